$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# FPVTL-2249 - add the "appointedGuardianClause" merge-field paragraph (and
# its surrounding blank spacer paragraphs) right before the existing
# "futureHearingClause" merge-field paragraph near the end of the document.
# ---------------------------------------------------------------------------

# Locate the "futureHearingClause" text robustly (rather than assuming a
# fixed paragraph index) using Find.
$findRange = $d.Content
$found = $findRange.Find.Execute("futureHearingClause", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Walk the Paragraphs collection to get the Paragraph object that contains
# the match (the paragraph with "{{futureHearingClause}}").
$paraCount = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $candidate = $d.Paragraphs.Item($i).Range
    if ($candidate.Start -le $findRange.Start -and $candidate.End -ge $findRange.End) {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate the 'futureHearingClause' paragraph."
}

$targetPara = $d.Paragraphs.Item($targetIndex)
$targetRange = $targetPara.Range

# Insert the new "{{appointedGuardianClause}}" paragraph plus a following
# blank paragraph immediately before the "futureHearingClause" paragraph,
# in one shot so the new text lands in its own run (matching the way the
# rest of the template's merge-field paragraphs are authored).
$insertionPoint = $targetRange.Duplicate
$insertionPoint.Collapse(1)
$insertionPoint.InsertBefore("{{appointedGuardianClause}}" + [char]13 + [char]13)

# Re-resolve the (now shifted) "futureHearingClause" paragraph and append a
# blank paragraph immediately after it.
$paraCount2 = $d.Paragraphs.Count
for ($i = 1; $i -le $paraCount2; $i++) {
    $candidate2 = $d.Paragraphs.Item($i).Range
    if ($candidate2.Text -eq "{{futureHearingClause}}" + [char]13) {
        $fhcIndex = $i
        break
    }
}

if (-not $fhcIndex) {
    throw "Could not re-locate the 'futureHearingClause' paragraph after inserting the new text."
}

$fhcPara = $d.Paragraphs.Item($fhcIndex)
$fhcPara.Range.InsertParagraphAfter()

Write-Output "Inserted appointedGuardianClause paragraph + spacer paragraphs."
